$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
# for rows 2..21, columns D, J, K, L, M, P
$data = @{
  2  = @{ D = 44313; J = 250; K = 900;  L = 1000; M = 950;  P = 950 }
  3  = @{ D = 44249; J = 200; K = 900;  L = 1000; M = 950;  P = 950 }
  4  = @{ D = 44442; J = 240; K = 2300; L = 2500; M = 2400; P = 2400 }
  5  = @{ D = 44302; J = 200; K = 900;  L = 1000; M = 950;  P = 950 }
  6  = @{ D = 44365; J = 250; K = 2400; L = 2500; M = 2450; P = 2450 }
  7  = @{ D = 44376; J = 270; K = 2400; L = 2500; M = 2437; P = 2437 }
  8  = @{ D = 44260; J = 250; K = 900;  L = 1000; M = 950;  P = 950 }
  9  = @{ D = 44417; J = 250; K = 4000; L = 4500; M = 4250; P = 4250 }
  10 = @{ D = 44330; J = 250; K = 2800; L = 3000; M = 2900; P = 2900 }
  11 = @{ D = 44224; J = 200; K = 750;  L = 800;  M = 775;  P = 775 }
  12 = @{ D = 44292; J = 250; K = 1800; L = 2000; M = 1900; P = 1900 }
  13 = @{ D = 44280; J = 250; K = 1400; L = 1500; M = 1450; P = 1450 }
  14 = @{ D = 44274; J = 250; K = 1000; L = 1200; M = 1100; P = 1100 }
  15 = @{ D = 44474; J = 250; K = 2000; L = 2500; M = 2250; P = 2250 }
  16 = @{ D = 44362; J = 250; K = 2800; L = 3000; M = 2900; P = 2900 }
  17 = @{ D = 44349; J = 250; K = 2800; L = 3000; M = 2900; P = 2900 }
  18 = @{ D = 44435; J = 300; K = 2300; L = 2500; M = 2400; P = 2400 }
  19 = @{ D = 44250; J = 250; K = 1000; L = 1200; M = 1100; P = 1100 }
  20 = @{ D = 44326; J = 200; K = 2700; L = 2800; M = 2750; P = 2750 }
  21 = @{ D = 44432; J = 300; K = 2300; L = 2500; M = 2400; P = 2400 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("J$row").Value = $vals.J
  $ws.Range("K$row").Value = $vals.K
  $ws.Range("L$row").Value = $vals.L
  $ws.Range("M$row").Value = $vals.M
  $ws.Range("P$row").Value = $vals.P
}
